$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-30 down to 9-31.
$ws.Rows.Item(8).Insert()

# Copy formatting of the date cell (D9, formerly D8) onto the new D8 cell.
$ws.Range("D9").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 8 values.
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 45076
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100108
$ws.Cells.Item(8, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(8, 9).Value = 100108001
$ws.Cells.Item(8, 10).Value = "Guayaba"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 2600
$ws.Cells.Item(8, 15).Value = 2600
$ws.Cells.Item(8, 16).Value = 2600
$ws.Cells.Item(8, 17).Value = "$/kilo"
$ws.Cells.Item(8, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 19).Value = 2600
$ws.Cells.Item(8, 20).Value = 1
